$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.80"
$ws.Range("E2").Value = "'3.17%"
$ws.Range("E3").Value = "'2.64%"
$ws.Range("E4").Value = "'0.34%"
$ws.Range("D5").Value = "'0.08199"
$ws.Range("E5").Value = "'1.85%"
$ws.Range("D6").Value = "'1.961"
$ws.Range("E6").Value = "'1.44%"
$ws.Range("D7").Value = "'8.250"
$ws.Range("E7").Value = "'3.79%"
$ws.Range("D8").Value = "'0.9317"
$ws.Range("E8").Value = "'0.28%"
$ws.Range("D9").Value = "'0.1417"
$ws.Range("E9").Value = "'-1.76%"
$ws.Range("D10").Value = "'0.1973"
$ws.Range("E10").Value = "'2.84%"
$ws.Range("D11").Value = "'0.09095"
$ws.Range("E11").Value = "'1.23%"
$ws.Range("D12").Value = "'0.03544"
$ws.Range("E12").Value = "'0.45%"
$ws.Range("D13").Value = "'0.09813"
$ws.Range("E13").Value = "'0.45%"
$ws.Range("D14").Value = "'0.001404"
$ws.Range("E14").Value = "'0.58%"
$ws.Range("D15").Value = "'0.005973"
$ws.Range("E15").Value = "'-1.69%"
$ws.Range("E16").Value = "'-1.68%"
$ws.Range("D17").Value = "'4.271"
$ws.Range("E17").Value = "'1.74%"
$ws.Range("D18").Value = "'3.321"
$ws.Range("E18").Value = "'-3.04%"
$ws.Range("E19").Value = "'0.00%"
$ws.Range("D20").Value = "'0.1302"
$ws.Range("E20").Value = "'-2.40%"
$ws.Range("D21").Value = "'4.873"
$ws.Range("E21").Value = "'0.79%"
$ws.Range("D23").Value = "'0.04322"
$ws.Range("E23").Value = "'-0.68%"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'-0.65%"
$ws.Range("D25").Value = "'0.004788"
$ws.Range("E25").Value = "'16.33%"
$ws.Range("E26").Value = "'-0.32%"
$ws.Range("D27").Value = "'0.0003995"
$ws.Range("E27").Value = "'-10.18%"
$ws.Range("D39").Value = "'0.02244"
$ws.Range("E39").Value = "'7.86%"
$ws.Range("D40").Value = "'0.05281"
$ws.Range("E40").Value = "'4.89%"
$ws.Range("D41").Value = "'0.007549"
$ws.Range("E41").Value = "'1.17%"
$ws.Range("D42").Value = "'0.009861"
$ws.Range("E42").Value = "'-2.70%"
$ws.Range("D43").Value = "'0.1377"
$ws.Range("E43").Value = "'2.25%"
$ws.Range("D44").Value = "'0.002115"
$ws.Range("E44").Value = "'-1.33%"
$ws.Range("D45").Value = "'0.009793"
$ws.Range("E45").Value = "'8.41%"
$ws.Range("D46").Value = "'0.00006363"
$ws.Range("E46").Value = "'2.45%"
$ws.Range("E47").Value = "'-0.31%"
$ws.Range("E48").Value = "'-7.42%"
$ws.Range("E49").Value = "'-25.23%"
$ws.Range("E50").Value = "'-0.31%"
$ws.Range("E51").Value = "'-0.31%"
